# Update "SUPORT IMAGEN" workbook:
#  - Clear the "source" (column K) tag for rows that were tagged "Últimas
#    Unidades" (they no longer carry that label).
#  - Collapse the three different "Preventa 26/10|26/11|26/12" pre-sale
#    labels (rows 14-16) onto the same "Preventa 26/09" label already used
#    by row 13.
#  - Re-enter the column O helper formula across O2:O24 as a single fill
#    so it recalculates cleanly for every row.
#  - Make "Hoja2" (the data sheet) the active/selected sheet instead of
#    "Hoja1".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja2")

# Rows whose "source" tag ("Últimas Unidades") is being removed entirely.
$clearRows = @(6,7,8,9,10,11,12,17,18,19,20,21,22,23,24)
foreach ($r in $clearRows) {
    $ws.Range("K$r").Value = ""
}

# Rows 14-16 move from "Preventa 26/10"/"26/11"/"26/12" to "Preventa 26/09".
$ws.Range("K14").Value = "Preventa 26/09"
$ws.Range("K15").Value = "Preventa 26/09"
$ws.Range("K16").Value = "Preventa 26/09"

# Rebuild the concatenated CSV-row helper column as one fill so it becomes
# a single shared formula again.
$ws.Range("O2:O24").Formula = '=CONCATENATE(A2,",",B2,",",C2,",",D2,",",E2,",",F2,",",G2,",",H2,",",I2,",",J2,",",K2,",",L2)'

# Hoja2 becomes the active sheet/tab (was Hoja1).
$ws.Activate()
$ws.Range("A1").Select()
